# Add slide titles for all slides: insert 13 new "Title and Content" slides
# after the existing title slide, each carrying just a title (content
# placeholder left empty), matching the new lecture outline.

$p = $ppt.ActivePresentation

# ppLayoutText-ish "Title and Content" custom layout is the 2nd layout
# on the slide master (index 2) - same as used by the other body slides.
$layoutIndex = 2

$titles = @(
    "Moore’s law and genomics",
    "Trivial parallelization",
    "Parallel make",
    "GPU",
    "Clusters",
    "Message-passing parallelization",
    "MPI",
    "CORBA/COM",
    "Web services",
    "Grids",
    "Condor",
    "Example: mycorrhiza pipeline",
    "Exercise: parallelize the make file"
)

$insertAt = 2
foreach ($title in $titles) {
    $slide = $p.Slides.Add($insertAt, $layoutIndex)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $title
    $insertAt = $insertAt + 1
}

Write-Host "Total slides:" $p.Slides.Count
